$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 5 new rows before existing row 429 (old rows 429-473 shift down to 434-478)
$ws.Range("A429:A433").EntireRow.Insert()

# Populate the 5 newly inserted rows with the new price-report entries.
# Columns A-J are identical across this block (same market/product/category).

$commonA = "8"
$commonB = "Terminal La Palmera de La Serena"
$commonC = "Coquimbo"
$commonE = 4
$commonF = "Fruta"
$commonG = 100103
$commonH = "Frutos de hueso (carozo)"
$commonI = 100103004
$commonJ = "Durazno"

$newRows = @(
    @{Row=429; D=44918; K="Kurakata";   L="Especial"; M=20; N=495000; O=500000; P=497500; Q="`$/bins (420 kilos)"; R="Región de O'Higgins"; S=1185; T=420},
    @{Row=430; D=44918; K="Kurakata";   L="Primera";  M=16; N=460000; O=470000; P=465000; Q="`$/bins (420 kilos)"; R="Región de O'Higgins"; S=1107; T=420},
    @{Row=431; D=44918; K="Kurakata";   L="Segunda";  M=16; N=420000; O=430000; P=425000; Q="`$/bins (420 kilos)"; R="Región de O'Higgins"; S=1012; T=420},
    @{Row=432; D=44918; K="Toscana";    L="Especial"; M=16; N=490000; O=500000; P=495000; Q="`$/bins (420 kilos)"; R="Región de O'Higgins"; S=1179; T=420},
    @{Row=433; D=44918; K="Toscana";    L="Primera";  M=10; N=420000; O=430000; P=425000; Q="`$/bins (420 kilos)"; R="Región de O'Higgins"; S=1012; T=420}
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $commonA
    $ws.Cells.Item($row, 2).Value = $commonB
    $ws.Cells.Item($row, 3).Value = $commonC
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $commonE
    $ws.Cells.Item($row, 6).Value = $commonF
    $ws.Cells.Item($row, 7).Value = $commonG
    $ws.Cells.Item($row, 8).Value = $commonH
    $ws.Cells.Item($row, 9).Value = $commonI
    $ws.Cells.Item($row, 10).Value = $commonJ
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}
